# Generate Report for Handoff
# Updates the localization-status report workbook:
#  - Overview sheet: re-orders handoff rows and refreshes the status/timestamp
#    for the file that was just re-handed-off (10c11ab6-...) while the other
#    two files keep their "Handed back" status.
#  - zh-cn / de-de detail sheets: mirror the same re-ordering and refresh the
#    handoff detail row (cols A,B,C,D,E,F,G,H,J) for the refreshed file.
#  - Hyperlink display text is refreshed to track the new cell contents while
#    the underlying hyperlink targets (r:id / Address) are left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ffff328d9201-3368-4723-993e-08ce926f2040.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value = "2016-03-22 05:43:07"

$wsOverview.Range("A3").Value = "ffffffe9951488-836a-4f76-a359-4577551c8c41.md"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D3").Value = "2016-03-22 05:43:07"

$wsOverview.Range("A4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-03-22 05:48:21"

$overviewDisplays = @(
    "ffff328d9201-3368-4723-993e-08ce926f2040.md",
    "ffffffe9951488-836a-4f76-a359-4577551c8c41.md",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md"
)
$i = 0
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $overviewDisplays[$i]
    $i = $i + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "ffff328d9201-3368-4723-993e-08ce926f2040.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D2").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-22 05:42:59"
$wsZhCn.Range("F2").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md"
$wsZhCn.Range("G2").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-03-22 05:43:38"
$wsZhCn.Range("J2").Value = "Include"

$wsZhCn.Range("A3").Value = "ffffffe9951488-836a-4f76-a359-4577551c8c41.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-22 05:42:59"
$wsZhCn.Range("F3").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md"
$wsZhCn.Range("G3").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-03-22 05:43:38"
$wsZhCn.Range("J3").Value = "Include"

$wsZhCn.Range("A4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.ce8bd949dd05ab366c7ef93937f0d33a61034b1b.zh-cn.xlf"
$wsZhCn.Range("E4").Value = "2016-03-22 05:48:12"
$wsZhCn.Range("F4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md"
$wsZhCn.Range("G4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.ce8bd949dd05ab366c7ef93937f0d33a61034b1b.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-03-22 05:47:23"
$wsZhCn.Range("J4").Value = "Include"

$zhCnDisplays = @(
    "ffff328d9201-3368-4723-993e-08ce926f2040.md",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf",
    "ffffffe9951488-836a-4f76-a359-4577551c8c41.md",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.ce8bd949dd05ab366c7ef93937f0d33a61034b1b.zh-cn.xlf",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.ce8bd949dd05ab366c7ef93937f0d33a61034b1b.zh-cn.xlf"
)
$i = 0
foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = $zhCnDisplays[$i]
    $i = $i + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "ffff328d9201-3368-4723-993e-08ce926f2040.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D2").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-22 05:43:07"
$wsDeDe.Range("F2").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md"
$wsDeDe.Range("G2").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-03-22 05:43:52"
$wsDeDe.Range("J2").Value = "Include"

$wsDeDe.Range("A3").Value = "ffffffe9951488-836a-4f76-a359-4577551c8c41.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-22 05:43:07"
$wsDeDe.Range("F3").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md"
$wsDeDe.Range("G3").Value = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-03-22 05:43:52"
$wsDeDe.Range("J3").Value = "Include"

$wsDeDe.Range("A4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.ce8bd949dd05ab366c7ef93937f0d33a61034b1b.de-de.xlf"
$wsDeDe.Range("E4").Value = "2016-03-22 05:48:21"
$wsDeDe.Range("F4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md"
$wsDeDe.Range("G4").Value = "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.ce8bd949dd05ab366c7ef93937f0d33a61034b1b.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-03-22 05:47:36"
$wsDeDe.Range("J4").Value = "Include"

$deDeDisplays = @(
    "ffff328d9201-3368-4723-993e-08ce926f2040.md",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf",
    "ffffffe9951488-836a-4f76-a359-4577551c8c41.md",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md",
    "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.ce8bd949dd05ab366c7ef93937f0d33a61034b1b.de-de.xlf",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.md",
    "10c11ab6-ade2-44d1-bec5-e6f9c874c6dc.ce8bd949dd05ab366c7ef93937f0d33a61034b1b.de-de.xlf"
)
$i = 0
foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = $deDeDisplays[$i]
    $i = $i + 1
}

Write-Host "Report regenerated for handoff."
